$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data update (REPORT_DATE changed from 2019-12-31 to 2018-12-31,
# along with the corresponding cash-flow figures for that report date).
$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 2540490248.14
$ws.Range("P2").Value = 198.216789396
$ws.Range("Q2").Value = 43263468510.74
$ws.Range("R2").Value = 3375.547626138
$ws.Range("S2").Value = 1737859091.88
$ws.Range("T2").Value = 135.593061169
$ws.Range("U2").Value = -3158011726.17
$ws.Range("V2").Value = -246.3976965449
$ws.Range("Y2").Value = 3065054955.35
$ws.Range("Z2").Value = 239.1449260696
$ws.Range("AA2").Value = 1863009126.88
$ws.Range("AB2").Value = 145.3576482004
$ws.Range("AC2").Value = 1281672584.79
$ws.Range("AD2").Value = -3.7085237094
